$d = $word.ActiveDocument

$d.Content.Find.Execute("Roy Batty, August 17 2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lewis Comstive, July 2021", 2)
